# Update automatico via Actualizar 06-03-2020 02-45-10
# Adds a new daily record (2020-06-02) to the "Condicion_Pacientes" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the Excel table by one row; this keeps the table ref/autoFilter,
# the sheet dimension, etc. in sync automatically.
$lo = $ws.ListObjects.Item("Condicion_Pacientes")
$newRow = $lo.ListRows.Add()

# Copy the formatting (number formats / styles) from the previous last
# row (81) down into the freshly added row (82) before filling values.
$ws.Range("A81:F81").Copy() | Out-Null
$ws.Range("A82:F82").PasteSpecial(-4122) | Out-Null

# New data for 2020-06-02 (serial date 43984)
$ws.Range("A82").Value = 43984
$ws.Range("B82").Value = 573
$ws.Range("C82").Value = 165
$ws.Range("D82").Value = 408
$ws.Range("E82").Value = 98
$ws.Range("F82").Value = 44

# Match the saved selection state from the edit.
$ws.Range("D73").Select() | Out-Null
